# Regenerate the "K" column (column G) values for rows 2-77.
# The new values below were computed by the author's script (std/mean
# recalculation of s_vals) and replace the previous Strike# derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(
    0,1,0,1,2,1,0,2,0,0,
    0,2,3,2,0,3,1,3,0,0,
    0,1,1,1,1,1,0,2,1,0,
    1,0,2,0,0,0,2,1,0,1,
    1,3,1,2,2,0,0,0,0,1,
    0,0,3,0,2,2,1,0,0,0,
    1,2,0,0,1,1,1,1,1,1,
    2,0,3,1,1,1
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
